$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Header row 2 (A2): now wraps its text too (same border as before, style
# already used by B2/C2)
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# New content rows 3 & 4 - two new tasks about understanding the usage
# context via interviews with students / professors.
#
# New shared strings must be minted in this exact order so the sharedStrings
# table lines up: C3, C4, B3(&B4), A3, A4
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "Cuong, Konstantin, Simon"
$ws.Range("C4").Value = "Linda, Carola"
$ws.Range("B3").Value = "Interviews/ Fokusgruppen"
$ws.Range("B4").Value = "Interviews/ Fokusgruppen"
$ws.Range("A3").Value = "Verstehen & Festlegen des Nutzungskontexts: Studenten"
$ws.Range("A4").Value = "Verstehen & Festlegen des Nutzungskontexts: Professoren"

$ws.Range("D3").Value = 43434
$ws.Range("E3").Value = "7 Tage"
$ws.Range("F3").Value = 43441

$ws.Range("D4").Value = 43434
$ws.Range("E4").Value = "7 Tage"
$ws.Range("F4").Value = 43441

# --- formatting -------------------------------------------------------
# Build each new style exactly once on a "seed" cell, then fan it out with
# copy/paste-special so every matching cell shares the same cellXfs entry
# (setting NumberFormat/Borders cell-by-cell mints a brand new style each
# time, which would fragment the style table).

# Date cells (D3/F3/D4/F4): thin box border + existing date format (numFmtId 14)
$ws.Range("D3").Borders.LineStyle = 1
$ws.Range("D3").NumberFormat = "mm-dd-yy"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F4").PasteSpecial($xlPasteFormats) | Out-Null

# Duration cells (E3/E4): same style already used by E2 ("7 Tage")
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E4").PasteSpecial($xlPasteFormats) | Out-Null

# Task / assignee columns (A3/B3/A4/B4): thin box border + wrap text
$ws.Range("A3").Borders.LineStyle = 1
$ws.Range("A3").WrapText = $true
$ws.Range("A3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B4").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# Remaining empty rows (5-22): columns A & B now also carry the thin border +
# wrap-text style used by the new task rows above
# ---------------------------------------------------------------------------
$ws.Range("A5:B22").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Application.CutCopyMode = $false

# Rows 3 & 4 grew taller to fit the wrapped two-line text
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30

# ---------------------------------------------------------------------------
# Column widths - re-fit after the content/format changes
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 37
$ws.Columns.Item(2).ColumnWidth = 46.57
$ws.Columns.Item(3).ColumnWidth = 23.71
$ws.Columns.Item(4).ColumnWidth = 9.29
$ws.Columns.Item(5).ColumnWidth = 5.71
$ws.Columns.Item(6).ColumnWidth = 9.29

# ---------------------------------------------------------------------------
# Selection moved as a side-effect of the edit session
# ---------------------------------------------------------------------------
$ws.Range("J8").Select() | Out-Null
